$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain string/text assignments (safe from numeric auto-conversion)
$ws.Cells.Item(2, 4).Value = "30.257.12"
$ws.Cells.Item(2, 5).Value = "  +0.20%  "
$ws.Cells.Item(3, 4).Value = "1.862.39"
$ws.Cells.Item(3, 5).Value = "  -0.15%  "
$ws.Cells.Item(4, 5).Value = "  +0.01%  "
$ws.Cells.Item(5, 5).Value = "  +0.58%  "
$ws.Cells.Item(6, 5).Value = "  +0.05%  "
$ws.Cells.Item(7, 5).Value = "  +1.20%  "
$ws.Cells.Item(8, 5).Value = "  +2.32%  "
$ws.Cells.Item(9, 5).Value = "  +0.57%  "
$ws.Cells.Item(10, 5).Value = "  +2.84%  "
$ws.Cells.Item(11, 5).Value = "  +0.93%  "
$ws.Cells.Item(12, 5).Value = "  +0.41%  "
$ws.Cells.Item(13, 4).Value = "1.857.69"
$ws.Cells.Item(13, 5).Value = "  -0.52%  "
$ws.Cells.Item(14, 5).Value = "  +0.82%  "
$ws.Cells.Item(15, 5).Value = "  +1.20%  "
$ws.Cells.Item(16, 5).Value = "  -4.50%  "
$ws.Cells.Item(17, 4).Value = "30.250.77"
$ws.Cells.Item(17, 5).Value = "  +0.17%  "
$ws.Cells.Item(18, 5).Value = "  +8.18%  "
$ws.Cells.Item(19, 5).Value = "  +0.08%  "
$ws.Cells.Item(20, 5).Value = "  +3.72%  "
$ws.Cells.Item(21, 4).Value = "2.105.11"
$ws.Cells.Item(21, 5).Value = "  -0.25%  "
$ws.Cells.Item(22, 5).Value = "  +0.03%  "
$ws.Cells.Item(23, 5).Value = "  -4.31%  "
$ws.Cells.Item(24, 5).Value = "  +0.14%  "
$ws.Cells.Item(25, 5).Value = "  +1.71%  "
$ws.Cells.Item(26, 5).Value = "  -0.04%  "
$ws.Cells.Item(27, 5).Value = "  -0.99%  "
$ws.Cells.Item(28, 5).Value = "  +1.58%  "
$ws.Cells.Item(29, 5).Value = "  +1.52%  "
$ws.Cells.Item(30, 5).Value = "  +2.90%  "
$ws.Cells.Item(31, 5).Value = "  -1.74%  "
$ws.Cells.Item(32, 5).Value = "  -0.25%  "
$ws.Cells.Item(33, 5).Value = "  -1.89%  "
$ws.Cells.Item(34, 5).Value = "  +0.49%  "
$ws.Cells.Item(35, 5).Value = "  +1.17%  "
$ws.Cells.Item(36, 5).Value = "  -0.46%  "
$ws.Cells.Item(37, 5).Value = "  -0.67%  "
$ws.Cells.Item(38, 5).Value = "  +1.87%  "
$ws.Cells.Item(39, 5).Value = "  +3.52%  "
$ws.Cells.Item(40, 5).Value = "  +1.01%  "
$ws.Cells.Item(41, 5).Value = "  +0.32%  "
$ws.Cells.Item(42, 5).Value = "  +0.04%  "
$ws.Cells.Item(43, 2).Value = "TrustWalletToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(43, 5).Value = "  -0.22%  "
$ws.Cells.Item(44, 2).Value = "TheSandbox"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(44, 5).Value = "  +0.01%  "
$ws.Cells.Item(45, 2).Value = "PaxDollar"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(45, 5).Value = "  -0.04%  "
$ws.Cells.Item(46, 5).Value = "  -0.71%  "
$ws.Cells.Item(47, 5).Value = "  -0.30%  "
$ws.Cells.Item(48, 5).Value = "  +1.30%  "
$ws.Cells.Item(49, 5).Value = "  +0.65%  "
$ws.Cells.Item(50, 5).Value = "  +0.47%  "
$ws.Cells.Item(51, 5).Value = "  +0.69%  "

# Numeric-looking text values for column D: use a text formula then convert
# to a literal value via copy/paste-special so Excel keeps them as text
# (matching the original inline/shared-string cell type) instead of numbers.
$ws.Cells.Item(4, 4).Formula = '="1.001"'
$ws.Cells.Item(5, 4).Formula = '="236.43"'
$ws.Cells.Item(6, 4).Formula = '="1.001"'
$ws.Cells.Item(7, 4).Formula = '="0.4721"'
$ws.Cells.Item(8, 4).Formula = '="0.2897"'
$ws.Cells.Item(9, 4).Formula = '="0.06577"'
$ws.Cells.Item(10, 4).Formula = '="21.85"'
$ws.Cells.Item(12, 4).Formula = '="97.89"'
$ws.Cells.Item(15, 4).Formula = '="0.6813"'
$ws.Cells.Item(16, 4).Formula = '="267.48"'
$ws.Cells.Item(18, 4).Formula = '="13.71"'
$ws.Cells.Item(19, 4).Formula = '="1.002"'
$ws.Cells.Item(20, 4).Formula = '="0.000007547"'
$ws.Cells.Item(23, 4).Formula = '="5.282"'
$ws.Cells.Item(24, 4).Formula = '="6.173"'
$ws.Cells.Item(25, 4).Formula = '="167.79"'
$ws.Cells.Item(26, 4).Formula = '="9.196"'
$ws.Cells.Item(27, 4).Formula = '="18.93"'
$ws.Cells.Item(28, 4).Formula = '="1.952"'
$ws.Cells.Item(29, 4).Formula = '="1.394"'
$ws.Cells.Item(30, 4).Formula = '="0.09983"'
$ws.Cells.Item(31, 4).Formula = '="4.342"'
$ws.Cells.Item(32, 4).Formula = '="1.471"'
$ws.Cells.Item(33, 4).Formula = '="4.021"'
$ws.Cells.Item(34, 4).Formula = '="0.04710"'
$ws.Cells.Item(35, 4).Formula = '="1.131"'
$ws.Cells.Item(36, 4).Formula = '="0.7026"'
$ws.Cells.Item(37, 4).Formula = '="2.708"'
$ws.Cells.Item(38, 4).Formula = '="0.01882"'
$ws.Cells.Item(39, 4).Formula = '="2.626"'
$ws.Cells.Item(40, 4).Formula = '="6.320"'
$ws.Cells.Item(41, 4).Formula = '="74.01"'
$ws.Cells.Item(42, 4).Formula = '="1.945"'
$ws.Cells.Item(43, 4).Formula = '="0.8420"'
$ws.Cells.Item(44, 4).Formula = '="0.4169"'
$ws.Cells.Item(45, 4).Formula = '="1.000"'
$ws.Cells.Item(46, 4).Formula = '="103.30"'
$ws.Cells.Item(47, 4).Formula = '="7.159"'
$ws.Cells.Item(48, 4).Formula = '="946.26"'
$ws.Cells.Item(49, 4).Formula = '="9.196"'
$ws.Cells.Item(50, 4).Formula = '="34.18"'
$ws.Cells.Item(51, 4).Formula = '="0.05661"'

$rng = $ws.Range("D2:D51")
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = 0

